# Generate Report for Handback
#
# For each localized-file row on the Overview sheet and on each locale
# sheet (zh-cn, de-de), mark the file as handed back: update the status
# text and fill in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: columns E (zh-cn) and F (de-de) hold the per-locale
# status for each source file row.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# ---------------------------------------------------------------------
# Per-locale sheets: Status column (C), Latest Target File (I),
# Latest Handback File (J), Latest Handback DateTime (K).
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; HandbackDate = "2016-09-04 11:07:29" },
    @{ Name = "de-de"; HandbackDate = "2016-09-04 11:07:36" }
)

$rows = @(
    @{ Row = 2; Md = "e3bdac3a-024f-43dc-ae44-a88038eea087.md"; Xlf = "e3bdac3a-024f-43dc-ae44-a88038eea087.57c70ca50dcfb1dde6afab102524c8926f3165fb" },
    @{ Row = 3; Md = "fab465cc-28cd-467d-851c-5d24c32556e2.md"; Xlf = "fab465cc-28cd-467d-851c-5d24c32556e2.898b56de955ec9f412252b34d52ca8eee7320a1c" }
)

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53d8236ed1a406d357055437a167ac06e84b38b3/e2e/"

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    foreach ($r in $rows) {
        $row = $r.Row

        $ws.Range("C$row").Value = $newStatus
        $ws.Range("J$row").Value = "$($r.Xlf).$($locale.Name).xlf"
        $ws.Range("K$row").Value = $locale.HandbackDate

        $target = $ws.Range("I$row")
        $url = "$ghBase$($r.Md)"
        $ws.Hyperlinks.Add($target, $url, "", "", $r.Md) | Out-Null
    }
}
